$wb = $excel.ActiveWorkbook

# --- Sheet "positive" (sheet3.xml): add two new rows for amusement / relief ---
# (written first so the new shared strings land in the same order as upstream)
$wsPositive = $wb.Worksheets.Item("positive")

$wsPositive.Range("A15").Value = "amusement"
$wsPositive.Range("B15").Value = "en"
$wsPositive.Range("C15").Value = "amu"

$wsPositive.Range("A16").Value = "relief"
$wsPositive.Range("B16").Value = "en"
$wsPositive.Range("C16").Value = "rel"

$wsPositive.Activate()
$wsPositive.Range("B17").Select()

# --- Sheet "samples_retained" (sheet1.xml): fill in row 19 (LimaCastroScott) ---
$wsSamples = $wb.Worksheets.Item("samples_retained")

$wsSamples.Range("D19").Value = 59
$wsSamples.Range("E19").Value = 58
$wsSamples.Range("F19").Value = 0
$wsSamples.Range("H19").Formula = "=IF(OR(ISBLANK(D19), ISBLANK(E19),ISBLANK(F19)), " + [char]34 + [char]34 + ", SUM(D19:F19))"
$wsSamples.Range("I19").Value = 4
$wsSamples.Range("J19").Value = "anger, disgust, fear, sadness, achievement, amusement, pleasure, relief"
$wsSamples.Range("K19").Value = "4 pos 4 neg; discarded four for valence mismatch: relief_MS_13, relief_MS_14, relief_MS_15, fear_T_16"

# Restore the selection / view state on this sheet
$wsSamples.Activate()
$wsSamples.Range("D21").Select()

# --- Sheet "negative" (sheet4.xml): clear stored selection (select A1) ---
$wsNegative = $wb.Worksheets.Item("negative")
$wsNegative.Activate()
$wsNegative.Range("A1").Select()

# Reselect the samples_retained sheet as the active tab, matching tabSelected="1"
$wsSamples.Activate()
